$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 7: "STAMP Session-Sender Test Packet for SR Policy"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)

# Shape 4 = "Content Placeholder 2" - reposition / resize the text box
$shp7 = $s7.Shapes.Item(4)
$shp7.Left   = 304800 / 12700
$shp7.Top    = 1071168 / 12700
$shp7.Width  = 4478154 / 12700
$shp7.Height = 3879363 / 12700

$tr7 = $shp7.TextFrame.TextRange

# Paragraph 7: "Color only SR-MPLS Policy:" -> "Color-Only Destination Steering:"
$para7 = $tr7.Paragraphs(7, 1)
$para7.Text = "Color-Only Destination Steering:"

# Paragraph 8: "Destination Address in 127/8 address" -> "IPv4 "
$para8 = $tr7.Paragraphs(8, 1)
$para8.Text = "IPv4 "

# Insert two new level-2 paragraphs after paragraph 8
$newRange = $para8.InsertAfter("`rDestination Address in 127/8 range`rTTL is set to 1")

$para9 = $tr7.Paragraphs(9, 1)
$para9.IndentLevel = 3

$para10 = $tr7.Paragraphs(10, 1)
$para10.IndentLevel = 3

# Paragraph 11 (was "IPv4 TTL is set to 1") -> "IPv6 "
$para11 = $tr7.Paragraphs(11, 1)
$para11.Text = "IPv6 "

# Insert two more new level-2 paragraphs after paragraph 11
$newRange2 = $para11.InsertAfter("`rDestination Address set to ::1/128`rHop Limit is set to 1")

$para12 = $tr7.Paragraphs(12, 1)
$para12.IndentLevel = 3
$run12b = $para12.Runs(2, 1)
$run12b.Text = "::1/128"
$run12b.LanguageID = "en-CA"

$para13 = $tr7.Paragraphs(13, 1)
$para13.IndentLevel = 3
$run13a = $para13.Runs(1, 1)
$run13a.LanguageID = "en-CA"

# Split paragraph 12 run into "Destination Address set to " + "::1/128"
$run12a = $para12.Runs(1, 1)
$run12a.Text = "Destination Address set to "

# Uniformly resize every run in the text box from 14pt to 12pt
$tr7.Font.Size = 12

# ---------------------------------------------------------------------------
# Slide 8: "STAMP Session-Sender Test Packet for P2MP SR-MPLS Policy"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)

# Shape 1 = Title placeholder - shift vertical position slightly
$title8 = $s8.Shapes.Item(1)
$title8.Left   = 114300 / 12700
$title8.Top    = 102393 / 12700
$title8.Width  = 8915400 / 12700
$title8.Height = 845539 / 12700

# Shape 4 = "Content Placeholder 2" - update wording of the first paragraph
$shp8 = $s8.Shapes.Item(4)
$tr8 = $shp8.TextFrame.TextRange
$para8_1 = $tr8.Paragraphs(1, 1)
$para8_1.Text = "For end-to-end P2MP SR-MPLS Policy, an example STAMP Session-Sender test packet is sent with:"
